# Update "想去人数" (interest count) figures in the "展览" and "全部类型" sheets
# to match the newer scrape snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 382
$ws1.Range("F3").Value = 840
$ws1.Range("F4").Value = 285
$ws1.Range("F5").Value = 1045
$ws1.Range("F6").Value = 2432
$ws1.Range("F7").Value = 205

# --- Sheet "全部类型" (all types, combined listing) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 382
$ws4.Range("F3").Value = 840
$ws4.Range("F4").Value = 285
$ws4.Range("F7").Value = 1045
$ws4.Range("F8").Value = 2432
$ws4.Range("F10").Value = 205
